# Add a "percentage" column to both the "PI hours" and "dept hours" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("PI hours")
$ws2 = $wb.Worksheets.Item("dept hours")

# ---------------------------------------------------------------------
# Sheet "PI hours": insert a new column D ("percentage") before the
# existing "dept" column, which shifts "dept" from D to E.
# ---------------------------------------------------------------------
$ws1.Columns.Item(4).Insert()

# Copy header formatting (border/bold/alignment) from the "hours" header
# onto the new "percentage" header cell, then set its text.
$ws1.Range("C1").Copy()
$ws1.Range("D1").PasteSpecial(-4122)
$ws1.Range("D1").Value = "percentage"

$pi_total = 39 + 8 + 8
$ws1.Range("D2").Value = 39 / $pi_total * 100
$ws1.Range("D3").Value = 8 / $pi_total * 100
$ws1.Range("D4").Value = 8 / $pi_total * 100

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Sheet "dept hours": add a new column D ("percentage") after "hours".
# ---------------------------------------------------------------------
$ws2.Range("C1").Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("D1").Value = "percentage"

$dept_total = 55 + 39 + 39 + 8 + 8
$ws2.Range("D2").Value = 55 / $dept_total * 100
$ws2.Range("D3").Value = 39 / $dept_total * 100
$ws2.Range("D4").Value = 39 / $dept_total * 100
$ws2.Range("D5").Value = 8 / $dept_total * 100
$ws2.Range("D6").Value = 8 / $dept_total * 100

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Normalize the page margins on both sheets to Excel's standard defaults
# (0.7"/0.75"/0.3"), expressed in points since PageSetup margins are in
# points (1/72 inch).
# ---------------------------------------------------------------------
foreach ($ws in @($ws1, $ws2)) {
    $ws.PageSetup.LeftMargin = 50.4
    $ws.PageSetup.RightMargin = 50.4
    $ws.PageSetup.TopMargin = 54.0
    $ws.PageSetup.BottomMargin = 54.0
    $ws.PageSetup.HeaderMargin = 21.599999999999998
    $ws.PageSetup.FooterMargin = 21.599999999999998
}
